# The post at row 595 ("蟻に学べ") was removed from the source data.
# Delete that entire row; Excel will automatically shift all rows below it
# up by one (596 -> 595, 597 -> 596, ..., 794 -> 793) and the used range
# shrinks from A1:C794 to A1:C793.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(595).Delete()
